$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H98").Value = 2470.2942
$ws.Range("J98").Value = 3820
$ws.Range("L98").Value = 3820
$ws.Range("N98").Value = -6816

$ws.Range("H122").Value = 2470.2942
$ws.Range("J122").Value = 3820
$ws.Range("L122").Value = 11460
$ws.Range("N122").Value = -16360

$ws.Range("H125").Value = 1244.75
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 1326.3334
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 11937.0006
$ws.Range("M125").Value = -6540
$ws.Range("N125").Value = -16857.0006

$ws.Range("H132").Value = 6541388.5
$ws.Range("I132").Value = 9806253
$ws.Range("K132").Value = 29418759
$ws.Range("M132").Value = -29416229

$ws.Range("H138").Value = 691913.7
$ws.Range("J138").Value = 990209
$ws.Range("L138").Value = 2970627
$ws.Range("N138").Value = -2980907

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4135.294
$ws.Range("I32").Value = 4248.5806
$ws.Range("K32").Value = 4248.5806
$ws.Range("M32").Value = -3961.5806

$ws.Range("H88").Value = 2426.611
$ws.Range("I88").Value = 2008.7142
$ws.Range("J88").Value = 2692.5454
$ws.Range("K88").Value = 2008.7142
$ws.Range("L88").Value = 2692.5454
$ws.Range("M88").Value = -1602.7142
$ws.Range("N88").Value = -3504.5454

$ws.Range("H91").Value = 2426.611
$ws.Range("I91").Value = 2008.7142
$ws.Range("J91").Value = 2692.5454
$ws.Range("K91").Value = 2008.7142
$ws.Range("L91").Value = 2692.5454
$ws.Range("M91").Value = -604.7141999999999
$ws.Range("N91").Value = -5500.5454

$ws.Range("H97").Value = 995.1429000000001
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 1000
$ws.Range("N97").Value = -1992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1688.75
$ws.Range("I20").Value = 1650.8667
$ws.Range("J20").Value = 1751.8889
$ws.Range("K20").Value = 1650.8667
$ws.Range("L20").Value = 1751.8889
$ws.Range("M20").Value = -1403.8667
$ws.Range("N20").Value = -2245.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 62500876
$ws.Range("I16").Value = 71429390
$ws.Range("J16").Value = 1250
$ws.Range("K16").Value = 71429390
$ws.Range("L16").Value = 1250
$ws.Range("M16").Value = -71429103
$ws.Range("N16").Value = -1824

$ws.Range("H31").Value = 1593.5264
$ws.Range("I31").Value = 1473.5625
$ws.Range("J31").Value = 2233.3333
$ws.Range("K31").Value = 1473.5625
$ws.Range("L31").Value = 2233.3333
$ws.Range("M31").Value = -1178.5625
$ws.Range("N31").Value = -2823.3333

$ws.Range("H34").Value = 1593.5264
$ws.Range("I34").Value = 1473.5625
$ws.Range("J34").Value = 2233.3333
$ws.Range("K34").Value = 1473.5625
$ws.Range("L34").Value = 2233.3333
$ws.Range("M34").Value = -1271.5625
$ws.Range("N34").Value = -2637.3333

$ws.Range("H58").Value = 747.6875
$ws.Range("I58").Value = 760.6429000000001
$ws.Range("J58").Value = 657
$ws.Range("K58").Value = 760.6429000000001
$ws.Range("L58").Value = 657
$ws.Range("M58").Value = -557.6429000000001
$ws.Range("N58").Value = -1063

$ws.Range("H113").Value = 62500876
$ws.Range("I113").Value = 71429390
$ws.Range("J113").Value = 1250
$ws.Range("K113").Value = 71429390
$ws.Range("L113").Value = 1250
$ws.Range("M113").Value = -71427220
$ws.Range("N113").Value = -5590

$ws.Range("H132").Value = 5920.1924
$ws.Range("I132").Value = 7186.1177
$ws.Range("K132").Value = 21558.3531
$ws.Range("M132").Value = -19028.3531

$ws.Range("H134").Value = 14494332
$ws.Range("I134").Value = 17545296
$ws.Range("J134").Value = 2250
$ws.Range("K134").Value = 52635888
$ws.Range("L134").Value = 6750
$ws.Range("M134").Value = -52633353
$ws.Range("N134").Value = -11820

$ws.Range("H136").Value = 747.6875
$ws.Range("I136").Value = 760.6429000000001
$ws.Range("J136").Value = 657
$ws.Range("K136").Value = 2281.9287
$ws.Range("L136").Value = 1971
$ws.Range("M136").Value = 268.0712999999996
$ws.Range("N136").Value = -7071

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 3300
$ws.Range("J55").Value = 3300
$ws.Range("L55").Value = 9900
$ws.Range("N55").Value = -10254

$ws.Range("H69").Value = 2415.2856
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 2415.2856
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 7245.8568
$ws.Range("N69").Value = -8867.856800000001
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 2415.2856
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 2415.2856
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 21737.5704
$ws.Range("N72").Value = -29849.5704
$ws.Range("M72").ClearContents()

$ws.Range("H131").Value = 922
$ws.Range("I131").Value = 408.14285
$ws.Range("J131").Value = 960.6774
$ws.Range("K131").Value = 1224.42855
$ws.Range("L131").Value = 2882.0322
$ws.Range("M131").Value = 3815.57145
$ws.Range("N131").Value = -12962.0322

$ws.Range("H139").Value = 3000.5334
$ws.Range("I139").Value = 3054.4614
$ws.Range("J139").Value = 2650
$ws.Range("K139").Value = 9163.3842
$ws.Range("L139").Value = 7950
$ws.Range("M139").Value = -4023.3842
$ws.Range("N139").Value = -18230

$ws.Range("H141").Value = 3392.1667
$ws.Range("I141").Value = 1784.3334
$ws.Range("K141").Value = 5353.0002
$ws.Range("M141").Value = -173.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2951
$ws.Range("I132").Value = 3289
$ws.Range("J132").Value = 2477.8
$ws.Range("K132").Value = 9867
$ws.Range("L132").Value = 7433.400000000001
$ws.Range("M132").Value = -7337
$ws.Range("N132").Value = -12493.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 85000
$ws.Range("J69").Value = 85000
$ws.Range("L69").Value = 85000
$ws.Range("N69").Value = -86622

$ws.Range("H72").Value = 85000
$ws.Range("J72").Value = 85000
$ws.Range("L72").Value = 255000
$ws.Range("N72").Value = -263112

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H141").Value = 50165
$ws.Range("J141").Value = 50165
$ws.Range("L141").Value = 50165
$ws.Range("N141").Value = -60525

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 603.6667
$ws.Range("I136").Value = 304.4
$ws.Range("J136").Value = 2100
$ws.Range("K136").Value = 913.1999999999999
$ws.Range("L136").Value = 6300
$ws.Range("M136").Value = 1636.8
$ws.Range("N136").Value = -11400
